# Add 2022-Q1 data
#
# The original workbook has a "总计" (Total) summary sheet as the 6th sheet.
# This edit repurposes that sheet to hold the new "2022-Q1" fund-detail data
# (matching the pattern used by the other quarterly sheets), and appends a
# brand-new "总计" sheet after it containing the refreshed summary table
# (the old summary rows plus a new 2022-Q1 row at the top).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Repurpose the existing "总计" sheet (position 6) into "2022-Q1"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(6)
$q1.Cells.Clear()
$q1.Name = "2022-Q1"

# Grab a cell that already carries the "header" style (bold, bordered,
# centered) used throughout the workbook so we can copy that exact format
# instead of inventing a new style.
$styleSource = $wb.Worksheets.Item(5).Range("B1")
$indexStyleSource = $wb.Worksheets.Item(5).Range("A2")

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$styleSource.Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Data rows -- codes and most numeric-looking figures are stored as text
# (matching the rest of the workbook), so force a text number format
# before assigning the values to stop Excel from coercing them to numbers.
$textCols = $q1.Range("B2:G4")
$textCols.NumberFormat = "@"

$q1.Range("B2").Value = "009447"
$q1.Range("C2").Value = "财通资管科技创新一年定期开放混合"
$q1.Range("D2").Value = "13.37"
$q1.Range("E2").Value = "88.02"
$q1.Range("F2").Value = "4.09"
$q1.Range("G2").Value = "0.5468"
$q1.Range("H2").Value = 4

$q1.Range("B3").Value = "013345"
$q1.Range("C3").Value = "富荣信息技术混合A"
$q1.Range("D3").Value = "1.96"
$q1.Range("E3").Value = "90.39"
$q1.Range("F3").Value = "4.49"
$q1.Range("G3").Value = "0.0880"
$q1.Range("H3").Value = 5

$q1.Range("B4").Value = "013346"
$q1.Range("C4").Value = "富荣信息技术混合C"
$q1.Range("D4").Value = "1.44"
$q1.Range("E4").Value = "90.39"
$q1.Range("F4").Value = "4.49"
$q1.Range("G4").Value = "0.0647"
$q1.Range("H4").Value = 5

# Column A holds the original row-index numbers, styled like the other
# sheets (bold/centered/bordered).
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2

$indexStyleSource.Copy()
$q1.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Add a brand new "总计" sheet after "2022-Q1" with the refreshed
#    summary table (2022-Q1 row added on top of the previous data).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$styleSource.Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.7

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 16
$total.Range("D3").Value = 3.03

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 1.52

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 0.28

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.23

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2020-Q4"
$total.Range("C7").Value = 2
$total.Range("D7").Value = 0.05

$indexStyleSource.Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$q1.Range("A1").Select()
